$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Overview sheet: zh-cn / de-de status columns (E2, F2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes (narrower "Status" columns) ---
# Target raw OOXML width is 13.4101845877511. Excel's ColumnWidth property is
# expressed in characters and is internally snapped to whole pixels
# (pixels = round(ColumnWidth * 6) + 5, stored width = pixels / 6), so we feed
# it the inverse of that so the serialized width lands on the nearest
# representable value (13.333333333333334, i.e. 80 px).
$newWidth = 13.4101845877511 - 5/6

# Overview sheet: columns E and F (zh-cn / de-de)
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# zh-cn sheet: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# de-de sheet: column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
